# Aktualizacja dokumentacji podłączenia kanałów analogowych i cyfrowych
# w dokumentacji (karta 6210) - arkusz "Arkusz1".
#
# The row 7 "Chan x IN" / "Digital x IN (/16)" pairs are shifted one slot
# to the left (the obsolete "Ohm IN (0V - 3,5V L - h)" reading is dropped),
# and the "Digital n IN (/16)" labels are regrouped together at the end of
# the row (AE7:AH7) instead of being interleaved with the "Chan n IN" cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Row 7: analog/digital channel relabeling ------------------------------

# Columns P7/Q7 ("Chan 1 IN" / "Digital 1 IN (/16)") are no longer used.
$ws.Range("P7").Clear()
$ws.Range("Q7").Clear()

# Shift "Chan n IN" labels one column group to the left.
$ws.Range("S7").Value = "Chan 1 IN"
$ws.Range("T7").Value = "Chan 2 IN"
$ws.Range("U7").Value = "Chan 3 IN"
$ws.Range("V7").Value = "Chan 4 IN"
$ws.Range("W7").Value = "Chan 5 IN"
$ws.Range("X7").Value = "Chan 6 IN"
$ws.Range("Z7").Value = "Chan 7 IN"
$ws.Range("AA7").Value = "Chan 8 IN"

# AB7 ("Chan 6 IN") is no longer used now that Chan 6 IN lives in X7.
$ws.Range("AB7").Clear()

# Regroup the "Digital n IN (/16)" labels together at the end of the row.
$ws.Range("AE7").Value = "Digital 2 IN (/16)"
$ws.Range("AF7").Value = "Digital 3 IN (/16)"
$ws.Range("AG7").Value = "Digital 1 IN (/16)"
$ws.Range("AH7").Value = "Digital 4 IN (/16)"

# --- Cosmetic: move the active selection to E7 ------------------------------
$ws.Range("E7").Select() | Out-Null
